$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @{ B = 0.04612272632165332;  C = 0.6122205438342831;  D = 0.7643369223465524;  E = 0.8742636457880153;  F = 0.8983570984052854;  G = 18 }
    3  = @{ B = -0.05484937042860891; C = 0.5505006181463401;  D = 0.631985967850102;   E = 0.7949754510990273;  F = 0.8174892103658266;  G = 17 }
    4  = @{ B = 0.0454318984618024;   C = 0.4533283699478702;  D = 0.4371073281485932;  E = 0.6611409291131454;  F = 0.6812093330742512;  G = 16 }
    5  = @{ B = 0.1619522260889591;   C = 0.3783700883078635;  D = 0.2480121053477303;  E = 0.4980081378328373;  F = 0.4874682207362279;  G = 15 }
    6  = @{ B = 0.1471867289230871;   C = 0.4060743056329635;  D = 0.2139495689888607;  E = 0.4625468289685496;  F = 0.455056961899273;   G = 14 }
    7  = @{ B = 0.169821430505846;    C = 0.3504591757809973;  D = 0.229570748894733;   E = 0.4791354181176058;  F = 0.4663250474243488;  G = 13 }
    8  = @{ B = 0.1844780524147233;   C = 0.3492283006690807;  D = 0.2072969125451589;  E = 0.4552987069443081;  F = 0.4347598217479242;  G = 12 }
    9  = @{ B = 0.2035825112292451;   C = 0.3072747587892817;  D = 0.2238621271029474;  E = 0.4731407053963413;  F = 0.4479485651801955;  G = 11 }
    10 = @{ B = 0.1565263708236955;   C = 0.2762321533226782;  D = 0.1325312503702182;  E = 0.3640484176180666;  F = 0.3464594662951709;  G = 10 }
    11 = @{ B = 0.1105734048294803;   C = 0.3275689177650167;  D = 0.2709241877029764;  E = 0.5205037826019869;  F = 0.5394765273654789;  G = 9 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    $ws.Range("B$row").Value = $rowData.B
    $ws.Range("C$row").Value = $rowData.C
    $ws.Range("D$row").Value = $rowData.D
    $ws.Range("E$row").Value = $rowData.E
    $ws.Range("F$row").Value = $rowData.F
    $ws.Range("G$row").Value = $rowData.G
}
